$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$lo = $ws.ListObjects.Item(1)

# Step 1: append a new (11th) table column cleanly at the end (K)
$lo.ListColumns.Add() | Out-Null
Write-Host "After Add():" $lo.Range.Address()

# Step 2: shift data in columns G:J (7:10) right into H:K (8:11), to make room at G
$src = $ws.Range($ws.Cells.Item(1,7), $ws.Cells.Item(60,10))
$dst = $ws.Range($ws.Cells.Item(1,8), $ws.Cells.Item(60,11))
$dst.Value = $src.Value()

# Step 3: clear column G (old data moved away)
$ws.Range($ws.Cells.Item(1,7), $ws.Cells.Item(60,7)).ClearContents()

Write-Host "=== headers now ==="
for ($c=1; $c -le 11; $c++) {
    Write-Host "col$c =" $ws.Cells.Item(1,$c).Value()
}

# Step 4: Set the new headers
$ws.Cells.Item(1,7).Value = "data_file_distances"
$ws.Cells.Item(1,8).Value = "data_file_zones"

Write-Host "=== ListColumns after header rename ==="
for ($i=1; $i -le $lo.ListColumns.Count; $i++) {
    Write-Host "ListColumn $i Name=" $lo.ListColumns.Item($i).Name " RangeAddr=" $lo.ListColumns.Item($i).Range.Address()
}
